# Apply the cell value updates described in the commit diff for Alexander_Profits
# (workbook sheets ALC, ARM, BSM, CUL, LTW, WVR; each row holds price/profit
# simulation data recalculated by the scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 250003500
$ws.Range("I32").Value = 333337340
$ws.Range("J32").Value = 2002
$ws.Range("K32").Value = 333337340
$ws.Range("L32").Value = 2002
$ws.Range("M32").Value = -333337014
$ws.Range("N32").Value = -2654

$ws.Range("H98").Value = 1865.8182
$ws.Range("I98").Value = 1055.2106
$ws.Range("J98").Value = 6999.6665
$ws.Range("K98").Value = 1055.2106
$ws.Range("L98").Value = 6999.6665
$ws.Range("M98").Value = 442.7893999999999
$ws.Range("N98").Value = -9995.666499999999

$ws.Range("H122").Value = 1865.8182
$ws.Range("I122").Value = 1055.2106
$ws.Range("J122").Value = 6999.6665
$ws.Range("K122").Value = 3165.6318
$ws.Range("L122").Value = 20998.9995
$ws.Range("M122").Value = -715.6318000000001
$ws.Range("N122").Value = -25898.9995

$ws.Range("H137").Value = 2501048.2
$ws.Range("I137").Value = 1112133.5
$ws.Range("J137").Value = 6667792
$ws.Range("K137").Value = 3336400.5
$ws.Range("L137").Value = 20003376
$ws.Range("M137").Value = -3333850.5
$ws.Range("N137").Value = -20008476

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2914018.8
$ws.Range("I32").Value = 4122.7295
$ws.Range("J32").Value = 20858378
$ws.Range("K32").Value = 4122.7295
$ws.Range("L32").Value = 20858378
$ws.Range("M32").Value = -3835.7295
$ws.Range("N32").Value = -20858952

$ws.Range("H121:L121").ClearContents()
$ws.Range("N121").ClearContents()

$ws.Range("H122:N122").ClearContents()

$ws.Range("H123:L123").ClearContents()
$ws.Range("N123").ClearContents()

$ws.Range("H124:L124").ClearContents()
$ws.Range("N124").ClearContents()

$ws.Range("H125:L125").ClearContents()
$ws.Range("N125").ClearContents()

$ws.Range("H126:L126").ClearContents()

$ws.Range("H127:L127").ClearContents()

$ws.Range("H128:L128").ClearContents()
$ws.Range("N128").ClearContents()

$ws.Range("H129:L129").ClearContents()
$ws.Range("N129").ClearContents()

$ws.Range("H130:L130").ClearContents()
$ws.Range("N130").ClearContents()

$ws.Range("H131:L131").ClearContents()
$ws.Range("N131").ClearContents()

$ws.Range("H132:N132").ClearContents()

$ws.Range("H133:L133").ClearContents()
$ws.Range("N133").ClearContents()

$ws.Range("H134:L134").ClearContents()
$ws.Range("N134").ClearContents()

$ws.Range("H135:L135").ClearContents()
$ws.Range("N135").ClearContents()

$ws.Range("H137:N137").ClearContents()

$ws.Range("H138:L138").ClearContents()
$ws.Range("N138").ClearContents()

$ws.Range("H139:N139").ClearContents()

$ws.Range("H140:L140").ClearContents()
$ws.Range("N140").ClearContents()

$ws.Range("H141:L141").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 804.36365
$ws.Range("I94").Value = 811.3684
$ws.Range("J94").Value = 760
$ws.Range("K94").Value = 811.3684
$ws.Range("L94").Value = 760
$ws.Range("M94").Value = -360.3684
$ws.Range("N94").Value = -1662

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 868.4211
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1375
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 4125
$ws.Range("M22").Value = -1331
$ws.Range("N22").Value = -4463

$ws.Range("H27").Value = 868.4211
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 1375
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 4125
$ws.Range("M27").Value = -1398
$ws.Range("N27").Value = -4329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5265161
$ws.Range("I7").Value = 8334763
$ws.Range("J7").Value = 2986.4285
$ws.Range("K7").Value = 8334763
$ws.Range("L7").Value = 2986.4285
$ws.Range("M7").Value = -8334651
$ws.Range("N7").Value = -3210.4285

$ws.Range("H100").Value = 2935.4375
$ws.Range("I100").Value = 3121.9167
$ws.Range("J100").Value = 2376
$ws.Range("K100").Value = 3121.9167
$ws.Range("L100").Value = 2376
$ws.Range("M100").Value = -2580.9167
$ws.Range("N100").Value = -3458

$ws.Range("H124:L124").ClearContents()

$ws.Range("H125:L125").ClearContents()
$ws.Range("N125").ClearContents()

$ws.Range("H126").Value = 5265161
$ws.Range("I126").Value = 8334763
$ws.Range("J126").Value = 2986.4285
$ws.Range("K126").Value = 25004289
$ws.Range("L126").Value = 8959.2855
$ws.Range("M126").Value = -25001819
$ws.Range("N126").Value = -13899.2855

$ws.Range("H127:L127").ClearContents()

$ws.Range("H128:L128").ClearContents()
$ws.Range("N128").ClearContents()

$ws.Range("H129:L129").ClearContents()
$ws.Range("N129").ClearContents()

$ws.Range("H130:L130").ClearContents()

$ws.Range("H131:L131").ClearContents()
$ws.Range("N131").ClearContents()

$ws.Range("H132:N132").ClearContents()

$ws.Range("H133:L133").ClearContents()
$ws.Range("N133").ClearContents()

$ws.Range("H134:L134").ClearContents()
$ws.Range("N134").ClearContents()

$ws.Range("H135:L135").ClearContents()
$ws.Range("N135").ClearContents()

$ws.Range("H136:N136").ClearContents()

$ws.Range("H137:L137").ClearContents()
$ws.Range("N137").ClearContents()

$ws.Range("H138:L138").ClearContents()
$ws.Range("N138").ClearContents()

$ws.Range("H139:L139").ClearContents()
$ws.Range("N139").ClearContents()

$ws.Range("H140:L140").ClearContents()
$ws.Range("N140").ClearContents()

$ws.Range("H141:L141").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 28549.166
$ws.Range("I113").Value = 43818.043
$ws.Range("J113").Value = 1535
$ws.Range("K113").Value = 131454.129
$ws.Range("L113").Value = 4605
$ws.Range("M113").Value = -129284.129
$ws.Range("N113").Value = -8945

$ws.Range("H136").Value = 1531.5
$ws.Range("I136").Value = 1555.5652
$ws.Range("J136").Value = 1393.125
$ws.Range("K136").Value = 4666.6956
$ws.Range("L136").Value = 4179.375
$ws.Range("M136").Value = -2116.6956
$ws.Range("N136").Value = -9279.375
